# endpoint se zona do imovel aceita tipo empreendimento HIS/HMP
#
# - rename sheet "categoria admitida pro zona" -> "categoria_admitida_por_zona"
# - add two new header columns on that sheet: A1 = "codigo_siszon", B1 = "sigla_zona"
#   (pushes the existing header columns one slot to the right, introducing two
#   new shared strings)
# - refresh the saved cursor/selection on every sheet, ending with
#   "categoria_admitida_por_zona" as the active tab

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Rename the 4th sheet
$ws4.Name = "categoria_admitida_por_zona"

# New header cells on that sheet (shifts old headers right, adds 2 shared strings)
$ws4.Range("A1").Value = "codigo_siszon"
$ws4.Range("B1").Value = "sigla_zona"

# Update the remembered selection on each sheet, in tab order
[void]$ws1.Activate()
[void]$ws1.Range("E20").Select()

[void]$ws2.Activate()
[void]$ws2.Range("G3").Select()

[void]$ws3.Activate()
[void]$ws3.Range("A2:M27").Select()

# Leave "categoria_admitida_por_zona" as the active sheet/tab
[void]$ws4.Activate()
[void]$ws4.Range("C6").Select()
